# Generate Report for Handoff
#
# A fresh handoff-report generation run refreshed the "Latest Handoff
# Datetime" timestamp for every file whose status is "Ready for handoff",
# and stamped their "Priority" column with the handoff type ("ht").
#
# Rows 7, 8, 9, 11, 12, 14 are the "Ready for handoff" rows that share the
# same (now stale) handoff timestamp on each localized-language sheet; rows
# 10 and 13 already have a different (still current) timestamp and are left
# untouched.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 11, 12, 14)

# "Overview" sheet - Latest HO Xliff Generate Date (column G) mirrors the
# de-de handoff datetime for these rows.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-27 22:22:33"
}

# "zh-cn" sheet - Latest Handoff Datetime (H) + Priority (E)
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Range("H$r").Value = "2016-08-27 22:22:27"
    $wsZh.Range("E$r").Value = "ht"
}

# "de-de" sheet - Latest Handoff Datetime (H) + Priority (E)
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Range("H$r").Value = "2016-08-27 22:22:33"
    $wsDe.Range("E$r").Value = "ht"
}
